# Auto-generated Excel COM-interop script applying scheduled market-data refresh
# to the Hyperion_Profits leve-profit tracker workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 39 (Leve Item ID 4603)
$ws.Range("H39").Value = 182.80952
$ws.Range("I39").Value = 68.59999999999999
$ws.Range("J39").Value = 286.63635
$ws.Range("K39").Value = 205.8
$ws.Range("L39").Value = 859.90905
$ws.Range("M39").Value = 90.20000000000002
$ws.Range("N39").Value = -1451.90905

# Row 80 (Leve Item ID 12605)
$ws.Range("H80").Value = 1304.48
$ws.Range("I80").Value = 890.55
$ws.Range("J80").Value = 2960.2
$ws.Range("K80").Value = 2671.65
$ws.Range("L80").Value = 8880.599999999999
$ws.Range("M80").Value = -1673.65
$ws.Range("N80").Value = -10876.6

# Row 83 (Leve Item ID 12605)
$ws.Range("H83").Value = 1304.48
$ws.Range("I83").Value = 890.55
$ws.Range("J83").Value = 2960.2
$ws.Range("K83").Value = 8014.95
$ws.Range("L83").Value = 26641.8
$ws.Range("M83").Value = -3022.95
$ws.Range("N83").Value = -36625.8

# Row 86 (Leve Item ID 12603)
$ws.Range("H86").Value = 1803.5428
$ws.Range("I86").Value = 2022.8
$ws.Range("J86").Value = 1511.2
$ws.Range("K86").Value = 2022.8
$ws.Range("L86").Value = 1511.2
$ws.Range("M86").Value = -899.8
$ws.Range("N86").Value = -3757.2

# Row 88 (Leve Item ID 12608)
$ws.Range("H88").Value = 2715.3225
$ws.Range("I88").Value = 1606
$ws.Range("J88").Value = 2928.6538
$ws.Range("K88").Value = 1606
$ws.Range("L88").Value = 2928.6538
$ws.Range("M88").Value = -1200
$ws.Range("N88").Value = -3740.6538

# Row 89 (Leve Item ID 12603)
$ws.Range("H89").Value = 1803.5428
$ws.Range("I89").Value = 2022.8
$ws.Range("J89").Value = 1511.2
$ws.Range("K89").Value = 10114
$ws.Range("L89").Value = 7556
$ws.Range("M89").Value = -4498
$ws.Range("N89").Value = -18788

# Row 91 (Leve Item ID 12608)
$ws.Range("H91").Value = 2715.3225
$ws.Range("I91").Value = 1606
$ws.Range("J91").Value = 2928.6538
$ws.Range("K91").Value = 1606
$ws.Range("L91").Value = 2928.6538
$ws.Range("M91").Value = -202
$ws.Range("N91").Value = -5736.6538

# Row 97 (Leve Item ID 19885)
$ws.Range("H97").Value = 1369
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1369
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 4107
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -5099

# Row 100 (Leve Item ID 19906)
$ws.Range("H100").Value = 2792.875
$ws.Range("I100").Value = 2792.875
$ws.Range("K100").Value = 2792.875
$ws.Range("M100").Value = -2251.875

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 3110.7295
$ws.Range("I138").Value = 2553.2354
$ws.Range("K138").Value = 7659.706200000001
$ws.Range("M138").Value = -2519.706200000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 4301.169
$ws.Range("I32").Value = 2711.673
$ws.Range("K32").Value = 2711.673
$ws.Range("M32").Value = -2424.673

# Row 63 (Leve Item ID 12528)
$ws.Range("H63").Value = 3939.4
$ws.Range("I63").Value = 3186.75
$ws.Range("K63").Value = 3186.75
$ws.Range("M63").Value = -2500.75

# Row 66 (Leve Item ID 12528)
$ws.Range("H66").Value = 3939.4
$ws.Range("I66").Value = 3186.75
$ws.Range("K66").Value = 15933.75
$ws.Range("M66").Value = -12501.75

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 2200
$ws.Range("I132").Value = 1348.6666
$ws.Range("J132").Value = 4434.75
$ws.Range("K132").Value = 4045.9998
$ws.Range("L132").Value = 13304.25
$ws.Range("M132").Value = -1515.9998
$ws.Range("N132").Value = -18364.25

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 4553149
$ws.Range("I86").Value = 5891294.5
$ws.Range("J86").Value = 3454
$ws.Range("K86").Value = 5891294.5
$ws.Range("L86").Value = 3454
$ws.Range("M86").Value = -5890171.5
$ws.Range("N86").Value = -5700

# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 4553149
$ws.Range("I89").Value = 5891294.5
$ws.Range("J89").Value = 3454
$ws.Range("K89").Value = 29456472.5
$ws.Range("L89").Value = 17270
$ws.Range("M89").Value = -29450856.5
$ws.Range("N89").Value = -28502

# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 3644646.8
$ws.Range("I94").Value = 6061572.5
$ws.Range("J94").Value = 19258.5
$ws.Range("K94").Value = 6061572.5
$ws.Range("L94").Value = 19258.5
$ws.Range("M94").Value = -6061121.5
$ws.Range("N94").Value = -20160.5

# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 3572502
$ws.Range("I107").Value = 4465363
$ws.Range("J107").Value = 1058.75
$ws.Range("K107").Value = 4465363
$ws.Range("L107").Value = 1058.75
$ws.Range("M107").Value = -4463443
$ws.Range("N107").Value = -4898.75

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 4986.087
$ws.Range("I134").Value = 2522.7058
$ws.Range("J134").Value = 11965.667
$ws.Range("K134").Value = 7568.117400000001
$ws.Range("L134").Value = 35897.001
$ws.Range("M134").Value = -5033.117400000001
$ws.Range("N134").Value = -40967.001

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 17939.822
$ws.Range("I31").Value = 1183.4
$ws.Range("J31").Value = 59830.875
$ws.Range("K31").Value = 1183.4
$ws.Range("L31").Value = 59830.875
$ws.Range("M31").Value = -888.4000000000001
$ws.Range("N31").Value = -60420.875

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 17939.822
$ws.Range("I34").Value = 1183.4
$ws.Range("J34").Value = 59830.875
$ws.Range("K34").Value = 1183.4
$ws.Range("L34").Value = 59830.875
$ws.Range("M34").Value = -981.4000000000001
$ws.Range("N34").Value = -60234.875

# Row 88 (Leve Item ID 10608)
$ws.Range("H88").Value = 34662
$ws.Range("I88").Value = 10999
$ws.Range("J88").Value = 46493.5
$ws.Range("K88").Value = 10999
$ws.Range("L88").Value = 46493.5
$ws.Range("M88").Value = -10593
$ws.Range("N88").Value = -47305.5

# Row 91 (Leve Item ID 10608)
$ws.Range("H91").Value = 34662
$ws.Range("I91").Value = 10999
$ws.Range("J91").Value = 46493.5
$ws.Range("K91").Value = 10999
$ws.Range("L91").Value = 46493.5
$ws.Range("M91").Value = -9595
$ws.Range("N91").Value = -49301.5

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 101686.266
$ws.Range("I132").Value = 68377.07000000001
$ws.Range("J132").Value = 226595.75
$ws.Range("K132").Value = 205131.21
$ws.Range("L132").Value = 679787.25
$ws.Range("M132").Value = -202601.21
$ws.Range("N132").Value = -684847.25

# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 3595.4
$ws.Range("I134").Value = 2853.2
$ws.Range("J134").Value = 5079.8
$ws.Range("K134").Value = 8559.599999999999
$ws.Range("L134").Value = 15239.4
$ws.Range("M134").Value = -6024.599999999999
$ws.Range("N134").Value = -20309.4

$ws = $wb.Worksheets.Item("CUL")
# Row 57 (Leve Item ID 4655)
$ws.Range("H57").Value = 6088.1113
$ws.Range("I57").Value = 2396.5
$ws.Range("J57").Value = 7142.857
$ws.Range("K57").Value = 7189.5
$ws.Range("L57").Value = 21428.571
$ws.Range("M57").Value = -6630.5
$ws.Range("N57").Value = -22546.571

# Row 69 (Leve Item ID 12850)
$ws.Range("H69").Value = 3266.3333
$ws.Range("J69").Value = 3394
$ws.Range("L69").Value = 10182
$ws.Range("N69").Value = -11804

# Row 72 (Leve Item ID 12850)
$ws.Range("H72").Value = 3266.3333
$ws.Range("J72").Value = 3394
$ws.Range("L72").Value = 30546
$ws.Range("N72").Value = -38658

# Row 133 (Leve Item ID 44073)
$ws.Range("H133").Value = 5999
$ws.Range("I133").Value = 5999
$ws.Range("K133").Value = 17997
$ws.Range("M133").Value = -12937

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 3489606.2
$ws.Range("I80").Value = 6099809.5
$ws.Range("J80").Value = 9335
$ws.Range("K80").Value = 6099809.5
$ws.Range("L80").Value = 9335
$ws.Range("M80").Value = -6098811.5
$ws.Range("N80").Value = -11331

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 3489606.2
$ws.Range("I83").Value = 6099809.5
$ws.Range("J83").Value = 9335
$ws.Range("K83").Value = 30499047.5
$ws.Range("L83").Value = 46675
$ws.Range("M83").Value = -30494055.5
$ws.Range("N83").Value = -56659

# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 426535.94
$ws.Range("I122").Value = 594650.5600000001
$ws.Range("J122").Value = 6249.5
$ws.Range("K122").Value = 1783951.68
$ws.Range("L122").Value = 18748.5
$ws.Range("M122").Value = -1781501.68
$ws.Range("N122").Value = -23648.5

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 9831.583000000001
$ws.Range("I40").Value = 7998
$ws.Range("K40").Value = 7998
$ws.Range("M40").Value = -7862

# Row 68 (Leve Item ID 12563)
$ws.Range("H68").Value = 3196.375
$ws.Range("I68").Value = 2714.4
$ws.Range("J68").Value = 3999.6667
$ws.Range("K68").Value = 2714.4
$ws.Range("L68").Value = 3999.6667
$ws.Range("M68").Value = -1965.4
$ws.Range("N68").Value = -5497.6667

# Row 71 (Leve Item ID 12563)
$ws.Range("H71").Value = 3196.375
$ws.Range("I71").Value = 2714.4
$ws.Range("J71").Value = 3999.6667
$ws.Range("K71").Value = 13572
$ws.Range("L71").Value = 19998.3335
$ws.Range("M71").Value = -9828
$ws.Range("N71").Value = -27486.3335

# Row 82 (Leve Item ID 12565)
$ws.Range("H82").Value = 5053236
$ws.Range("I82").Value = 6947600
$ws.Range("J82").Value = 1599
$ws.Range("K82").Value = 6947600
$ws.Range("L82").Value = 1599
$ws.Range("M82").Value = -6947239
$ws.Range("N82").Value = -2321

# Row 85 (Leve Item ID 12565)
$ws.Range("H85").Value = 5053236
$ws.Range("I85").Value = 6947600
$ws.Range("J85").Value = 1599
$ws.Range("K85").Value = 6947600
$ws.Range("L85").Value = 1599
$ws.Range("M85").Value = -6946352
$ws.Range("N85").Value = -4095

# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 83337464
$ws.Range("I93").Value = 83337464
$ws.Range("K93").Value = 83337464
$ws.Range("M93").Value = -83336216

# Row 107 (Leve Item ID 38752)
$ws.Range("H107").Value = 4874.5
$ws.Range("I107").Value = 4874.5
$ws.Range("K107").Value = 4874.5
$ws.Range("M107").Value = -2954.5

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 4222.6665
$ws.Range("I132").Value = 3452.575
$ws.Range("J132").Value = 8073.125
$ws.Range("K132").Value = 10357.725
$ws.Range("L132").Value = 24219.375
$ws.Range("M132").Value = -7827.724999999999
$ws.Range("N132").Value = -29279.375

# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 46595.652
$ws.Range("I136").Value = 59944.4
$ws.Range("J136").Value = 4122.364
$ws.Range("K136").Value = 179833.2
$ws.Range("L136").Value = 12367.092
$ws.Range("M136").Value = -177283.2
$ws.Range("N136").Value = -17467.092

$ws = $wb.Worksheets.Item("WVR")
# Row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 2622.0908
$ws.Range("J122").Value = 3599.8
$ws.Range("L122").Value = 10799.4
$ws.Range("N122").Value = -15699.4

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 2791.7693
$ws.Range("I136").Value = 2455.7576
$ws.Range("J136").Value = 4639.8335
$ws.Range("K136").Value = 7367.2728
$ws.Range("L136").Value = 13919.5005
$ws.Range("M136").Value = -4817.2728
$ws.Range("N136").Value = -19019.5005

